# "Colocando header nos graficos" - add a header label in column A (row 1)
# of each chart-data sheet, fix accented Portuguese text in the row labels
# (removing the bold/border "header" style from those labels, since it now
# only belongs to row 1 / the real header), update a couple of figures on
# the cost sheet, and drop the now-unused "Teto" row from the emissions
# sheet.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheets 1-4 share the same layout: a title row (B1:E1 = years) and a
# label column (A2:A12 = source/technology). Add "Fonte/Tecnologia" to
# A1 (styled like the rest of the header row) and clean up the labels.
# ---------------------------------------------------------------------
$sheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

$labels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New header cell for the label column, styled like the year headers.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial($xlPasteFormats)

    # Fix labels and drop the header-style formatting now that the bold
    # border belongs to row 1 only.
    foreach ($r in 2..12) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.Value = $labels[$r]
        $cell.ClearFormats()
    }
}

# ---------------------------------------------------------------------
# "Emissoes Totais (MtCO2eq)" sheet: add a header, fix labels, and
# remove the now-unused "Teto" row (row 4).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial($xlPasteFormats)

$ws5.Cells.Item(2, 1).Value = "P.Médio"
$ws5.Cells.Item(2, 1).ClearFormats()

$ws5.Cells.Item(3, 1).Value = "P.Crítico"
$ws5.Cells.Item(3, 1).ClearFormats()

$ws5.Rows("4:4").Delete()

# ---------------------------------------------------------------------
# "Custo Total (bilhões de R$)" sheet: add a header, relabel B1 as a
# year, fix the expansion-type labels and their updated figures.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial($xlPasteFormats)

# B1 needs to become the text "2015" (not the number 2015) while keeping
# its original bold/border style - force text via NumberFormat, then
# restore the clean style from a cell that already has it.
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws6.Range("A2").Copy()
$ws6.Range("B1").PasteSpecial($xlPasteFormats)

$ws6.Cells.Item(2, 1).Value = "Expansão Centralizada"
$ws6.Cells.Item(2, 1).ClearFormats()
$ws6.Cells.Item(2, 2).Value = 746

$ws6.Cells.Item(3, 1).Value = "Expansão por GD"
$ws6.Cells.Item(3, 1).ClearFormats()
$ws6.Cells.Item(3, 2).Value = 99
